$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13-131 down to 14-132
$ws.Rows("13:13").Insert()

# Populate the fixed (non-shifting) columns for the newly inserted row 13
$ws.Cells.Item(13, 1).Value = 8
$ws.Cells.Item(13, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 100112001
$ws.Cells.Item(13, 7).Value = "Berenjena"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# Set Fecha / Volumen / Precio min-max-prom / Unidad / Precio-Kg / Kg-o-Unidades
# for every data row 13-132 (row 13 is brand-new data; rows 14-132 take the value
# that used to sit one row above, completing the weekly-reorder shift described
# in the commit).
$ws.Cells.Item(13, 4).Value = 44670
$ws.Cells.Item(13, 10).Value = 500
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 9000
$ws.Cells.Item(13, 13).Value = 8500
$ws.Cells.Item(13, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(13, 16).Value = 170
$ws.Cells.Item(13, 17).Value = 50

$ws.Cells.Item(14, 4).Value = 44505
$ws.Cells.Item(14, 10).Value = 600
$ws.Cells.Item(14, 11).Value = 8000
$ws.Cells.Item(14, 12).Value = 9000
$ws.Cells.Item(14, 13).Value = 8500
$ws.Cells.Item(14, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(14, 16).Value = 142
$ws.Cells.Item(14, 17).Value = 60

$ws.Cells.Item(15, 4).Value = 44320
$ws.Cells.Item(15, 10).Value = 520
$ws.Cells.Item(15, 11).Value = 8000
$ws.Cells.Item(15, 12).Value = 9000
$ws.Cells.Item(15, 13).Value = 8500
$ws.Cells.Item(15, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(15, 16).Value = 142
$ws.Cells.Item(15, 17).Value = 60

$ws.Cells.Item(16, 4).Value = 44657
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 8500
$ws.Cells.Item(16, 12).Value = 9000
$ws.Cells.Item(16, 13).Value = 8750
$ws.Cells.Item(16, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(16, 16).Value = 175
$ws.Cells.Item(16, 17).Value = 50

$ws.Cells.Item(17, 4).Value = 44384
$ws.Cells.Item(17, 10).Value = 600
$ws.Cells.Item(17, 11).Value = 11000
$ws.Cells.Item(17, 12).Value = 12000
$ws.Cells.Item(17, 13).Value = 11500
$ws.Cells.Item(17, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(17, 16).Value = 192
$ws.Cells.Item(17, 17).Value = 60

$ws.Cells.Item(18, 4).Value = 44655
$ws.Cells.Item(18, 10).Value = 440
$ws.Cells.Item(18, 11).Value = 8000
$ws.Cells.Item(18, 12).Value = 9000
$ws.Cells.Item(18, 13).Value = 8500
$ws.Cells.Item(18, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(18, 16).Value = 170
$ws.Cells.Item(18, 17).Value = 50

$ws.Cells.Item(19, 4).Value = 44510
$ws.Cells.Item(19, 10).Value = 520
$ws.Cells.Item(19, 11).Value = 8000
$ws.Cells.Item(19, 12).Value = 8500
$ws.Cells.Item(19, 13).Value = 8250
$ws.Cells.Item(19, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(19, 16).Value = 138
$ws.Cells.Item(19, 17).Value = 60

$ws.Cells.Item(20, 4).Value = 44321
$ws.Cells.Item(20, 10).Value = 600
$ws.Cells.Item(20, 11).Value = 8000
$ws.Cells.Item(20, 12).Value = 9000
$ws.Cells.Item(20, 13).Value = 8500
$ws.Cells.Item(20, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(20, 16).Value = 142
$ws.Cells.Item(20, 17).Value = 60

$ws.Cells.Item(21, 4).Value = 44391
$ws.Cells.Item(21, 10).Value = 600
$ws.Cells.Item(21, 11).Value = 12000
$ws.Cells.Item(21, 12).Value = 13000
$ws.Cells.Item(21, 13).Value = 12500
$ws.Cells.Item(21, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(21, 16).Value = 208
$ws.Cells.Item(21, 17).Value = 60

$ws.Cells.Item(22, 4).Value = 44503
$ws.Cells.Item(22, 10).Value = 600
$ws.Cells.Item(22, 11).Value = 8000
$ws.Cells.Item(22, 12).Value = 8500
$ws.Cells.Item(22, 13).Value = 8250
$ws.Cells.Item(22, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(22, 16).Value = 138
$ws.Cells.Item(22, 17).Value = 60

$ws.Cells.Item(23, 4).Value = 44509
$ws.Cells.Item(23, 10).Value = 400
$ws.Cells.Item(23, 11).Value = 8000
$ws.Cells.Item(23, 12).Value = 9000
$ws.Cells.Item(23, 13).Value = 8500
$ws.Cells.Item(23, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(23, 16).Value = 142
$ws.Cells.Item(23, 17).Value = 60

$ws.Cells.Item(24, 4).Value = 44526
$ws.Cells.Item(24, 10).Value = 600
$ws.Cells.Item(24, 11).Value = 9000
$ws.Cells.Item(24, 12).Value = 10000
$ws.Cells.Item(24, 13).Value = 9500
$ws.Cells.Item(24, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(24, 16).Value = 158
$ws.Cells.Item(24, 17).Value = 60

$ws.Cells.Item(25, 4).Value = 44315
$ws.Cells.Item(25, 10).Value = 440
$ws.Cells.Item(25, 11).Value = 8000
$ws.Cells.Item(25, 12).Value = 9000
$ws.Cells.Item(25, 13).Value = 8500
$ws.Cells.Item(25, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(25, 16).Value = 142
$ws.Cells.Item(25, 17).Value = 60

$ws.Cells.Item(26, 4).Value = 44244
$ws.Cells.Item(26, 10).Value = 600
$ws.Cells.Item(26, 11).Value = 8000
$ws.Cells.Item(26, 12).Value = 9000
$ws.Cells.Item(26, 13).Value = 8500
$ws.Cells.Item(26, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(26, 16).Value = 142
$ws.Cells.Item(26, 17).Value = 60

$ws.Cells.Item(27, 4).Value = 44578
$ws.Cells.Item(27, 10).Value = 700
$ws.Cells.Item(27, 11).Value = 8000
$ws.Cells.Item(27, 12).Value = 9000
$ws.Cells.Item(27, 13).Value = 8500
$ws.Cells.Item(27, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(27, 16).Value = 170
$ws.Cells.Item(27, 17).Value = 50

$ws.Cells.Item(28, 4).Value = 44545
$ws.Cells.Item(28, 10).Value = 540
$ws.Cells.Item(28, 11).Value = 10000
$ws.Cells.Item(28, 12).Value = 11000
$ws.Cells.Item(28, 13).Value = 10500
$ws.Cells.Item(28, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(28, 16).Value = 175
$ws.Cells.Item(28, 17).Value = 60

$ws.Cells.Item(29, 4).Value = 44279
$ws.Cells.Item(29, 10).Value = 600
$ws.Cells.Item(29, 11).Value = 8000
$ws.Cells.Item(29, 12).Value = 9000
$ws.Cells.Item(29, 13).Value = 8500
$ws.Cells.Item(29, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(29, 16).Value = 142
$ws.Cells.Item(29, 17).Value = 60

$ws.Cells.Item(30, 4).Value = 44344
$ws.Cells.Item(30, 10).Value = 520
$ws.Cells.Item(30, 11).Value = 12000
$ws.Cells.Item(30, 12).Value = 13000
$ws.Cells.Item(30, 13).Value = 12500
$ws.Cells.Item(30, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(30, 16).Value = 208
$ws.Cells.Item(30, 17).Value = 60

$ws.Cells.Item(31, 4).Value = 44494
$ws.Cells.Item(31, 10).Value = 500
$ws.Cells.Item(31, 11).Value = 8000
$ws.Cells.Item(31, 12).Value = 9000
$ws.Cells.Item(31, 13).Value = 8500
$ws.Cells.Item(31, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(31, 16).Value = 142
$ws.Cells.Item(31, 17).Value = 60

$ws.Cells.Item(32, 4).Value = 44413
$ws.Cells.Item(32, 10).Value = 640
$ws.Cells.Item(32, 11).Value = 12000
$ws.Cells.Item(32, 12).Value = 13000
$ws.Cells.Item(32, 13).Value = 12500
$ws.Cells.Item(32, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(32, 16).Value = 208
$ws.Cells.Item(32, 17).Value = 60

$ws.Cells.Item(33, 4).Value = 44561
$ws.Cells.Item(33, 10).Value = 520
$ws.Cells.Item(33, 11).Value = 9000
$ws.Cells.Item(33, 12).Value = 10000
$ws.Cells.Item(33, 13).Value = 9500
$ws.Cells.Item(33, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(33, 16).Value = 158
$ws.Cells.Item(33, 17).Value = 60

$ws.Cells.Item(34, 4).Value = 44306
$ws.Cells.Item(34, 10).Value = 500
$ws.Cells.Item(34, 11).Value = 8500
$ws.Cells.Item(34, 12).Value = 9000
$ws.Cells.Item(34, 13).Value = 8750
$ws.Cells.Item(34, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(34, 16).Value = 146
$ws.Cells.Item(34, 17).Value = 60

$ws.Cells.Item(35, 4).Value = 44313
$ws.Cells.Item(35, 10).Value = 520
$ws.Cells.Item(35, 11).Value = 8000
$ws.Cells.Item(35, 12).Value = 9000
$ws.Cells.Item(35, 13).Value = 8500
$ws.Cells.Item(35, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(35, 16).Value = 142
$ws.Cells.Item(35, 17).Value = 60

$ws.Cells.Item(36, 4).Value = 44589
$ws.Cells.Item(36, 10).Value = 500
$ws.Cells.Item(36, 11).Value = 8000
$ws.Cells.Item(36, 12).Value = 9000
$ws.Cells.Item(36, 13).Value = 8500
$ws.Cells.Item(36, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(36, 16).Value = 170
$ws.Cells.Item(36, 17).Value = 50

$ws.Cells.Item(37, 4).Value = 44327
$ws.Cells.Item(37, 10).Value = 500
$ws.Cells.Item(37, 11).Value = 9000
$ws.Cells.Item(37, 12).Value = 10000
$ws.Cells.Item(37, 13).Value = 9500
$ws.Cells.Item(37, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(37, 16).Value = 158
$ws.Cells.Item(37, 17).Value = 60

$ws.Cells.Item(38, 4).Value = 44454
$ws.Cells.Item(38, 10).Value = 600
$ws.Cells.Item(38, 11).Value = 9000
$ws.Cells.Item(38, 12).Value = 10000
$ws.Cells.Item(38, 13).Value = 9500
$ws.Cells.Item(38, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(38, 16).Value = 158
$ws.Cells.Item(38, 17).Value = 60

$ws.Cells.Item(39, 4).Value = 44286
$ws.Cells.Item(39, 10).Value = 600
$ws.Cells.Item(39, 11).Value = 8000
$ws.Cells.Item(39, 12).Value = 9000
$ws.Cells.Item(39, 13).Value = 8500
$ws.Cells.Item(39, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(39, 16).Value = 142
$ws.Cells.Item(39, 17).Value = 60

$ws.Cells.Item(40, 4).Value = 44477
$ws.Cells.Item(40, 10).Value = 600
$ws.Cells.Item(40, 11).Value = 8000
$ws.Cells.Item(40, 12).Value = 9000
$ws.Cells.Item(40, 13).Value = 8500
$ws.Cells.Item(40, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(40, 16).Value = 142
$ws.Cells.Item(40, 17).Value = 60

$ws.Cells.Item(41, 4).Value = 44379
$ws.Cells.Item(41, 10).Value = 600
$ws.Cells.Item(41, 11).Value = 12000
$ws.Cells.Item(41, 12).Value = 13000
$ws.Cells.Item(41, 13).Value = 12500
$ws.Cells.Item(41, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(41, 16).Value = 208
$ws.Cells.Item(41, 17).Value = 60

$ws.Cells.Item(42, 4).Value = 44449
$ws.Cells.Item(42, 10).Value = 600
$ws.Cells.Item(42, 11).Value = 9000
$ws.Cells.Item(42, 12).Value = 10000
$ws.Cells.Item(42, 13).Value = 9500
$ws.Cells.Item(42, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(42, 16).Value = 158
$ws.Cells.Item(42, 17).Value = 60

$ws.Cells.Item(43, 4).Value = 44322
$ws.Cells.Item(43, 10).Value = 440
$ws.Cells.Item(43, 11).Value = 8000
$ws.Cells.Item(43, 12).Value = 9000
$ws.Cells.Item(43, 13).Value = 8500
$ws.Cells.Item(43, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(43, 16).Value = 142
$ws.Cells.Item(43, 17).Value = 60

$ws.Cells.Item(44, 4).Value = 44314
$ws.Cells.Item(44, 10).Value = 560
$ws.Cells.Item(44, 11).Value = 8000
$ws.Cells.Item(44, 12).Value = 9000
$ws.Cells.Item(44, 13).Value = 8500
$ws.Cells.Item(44, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(44, 16).Value = 142
$ws.Cells.Item(44, 17).Value = 60

$ws.Cells.Item(45, 4).Value = 44519
$ws.Cells.Item(45, 10).Value = 560
$ws.Cells.Item(45, 11).Value = 8000
$ws.Cells.Item(45, 12).Value = 8500
$ws.Cells.Item(45, 13).Value = 8250
$ws.Cells.Item(45, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(45, 16).Value = 138
$ws.Cells.Item(45, 17).Value = 60

$ws.Cells.Item(46, 4).Value = 44392
$ws.Cells.Item(46, 10).Value = 500
$ws.Cells.Item(46, 11).Value = 12000
$ws.Cells.Item(46, 12).Value = 13000
$ws.Cells.Item(46, 13).Value = 12500
$ws.Cells.Item(46, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(46, 16).Value = 208
$ws.Cells.Item(46, 17).Value = 60

$ws.Cells.Item(47, 4).Value = 44414
$ws.Cells.Item(47, 10).Value = 600
$ws.Cells.Item(47, 11).Value = 12500
$ws.Cells.Item(47, 12).Value = 13000
$ws.Cells.Item(47, 13).Value = 12750
$ws.Cells.Item(47, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(47, 16).Value = 212
$ws.Cells.Item(47, 17).Value = 60

$ws.Cells.Item(48, 4).Value = 44328
$ws.Cells.Item(48, 10).Value = 600
$ws.Cells.Item(48, 11).Value = 12000
$ws.Cells.Item(48, 12).Value = 13000
$ws.Cells.Item(48, 13).Value = 12500
$ws.Cells.Item(48, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(48, 16).Value = 208
$ws.Cells.Item(48, 17).Value = 60

$ws.Cells.Item(49, 4).Value = 44400
$ws.Cells.Item(49, 10).Value = 600
$ws.Cells.Item(49, 11).Value = 11500
$ws.Cells.Item(49, 12).Value = 12000
$ws.Cells.Item(49, 13).Value = 11750
$ws.Cells.Item(49, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(49, 16).Value = 196
$ws.Cells.Item(49, 17).Value = 60

$ws.Cells.Item(50, 4).Value = 44377
$ws.Cells.Item(50, 10).Value = 600
$ws.Cells.Item(50, 11).Value = 12000
$ws.Cells.Item(50, 12).Value = 13000
$ws.Cells.Item(50, 13).Value = 12500
$ws.Cells.Item(50, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(50, 16).Value = 208
$ws.Cells.Item(50, 17).Value = 60

$ws.Cells.Item(51, 4).Value = 44665
$ws.Cells.Item(51, 10).Value = 400
$ws.Cells.Item(51, 11).Value = 8000
$ws.Cells.Item(51, 12).Value = 9000
$ws.Cells.Item(51, 13).Value = 8500
$ws.Cells.Item(51, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(51, 16).Value = 170
$ws.Cells.Item(51, 17).Value = 50

$ws.Cells.Item(52, 4).Value = 44351
$ws.Cells.Item(52, 10).Value = 520
$ws.Cells.Item(52, 11).Value = 11500
$ws.Cells.Item(52, 12).Value = 12000
$ws.Cells.Item(52, 13).Value = 11750
$ws.Cells.Item(52, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(52, 16).Value = 196
$ws.Cells.Item(52, 17).Value = 60

$ws.Cells.Item(53, 4).Value = 44253
$ws.Cells.Item(53, 10).Value = 840
$ws.Cells.Item(53, 11).Value = 8000
$ws.Cells.Item(53, 12).Value = 8500
$ws.Cells.Item(53, 13).Value = 8250
$ws.Cells.Item(53, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(53, 16).Value = 138
$ws.Cells.Item(53, 17).Value = 60

$ws.Cells.Item(54, 4).Value = 44350
$ws.Cells.Item(54, 10).Value = 500
$ws.Cells.Item(54, 11).Value = 12000
$ws.Cells.Item(54, 12).Value = 12500
$ws.Cells.Item(54, 13).Value = 12250
$ws.Cells.Item(54, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(54, 16).Value = 204
$ws.Cells.Item(54, 17).Value = 60

$ws.Cells.Item(55, 4).Value = 44399
$ws.Cells.Item(55, 10).Value = 600
$ws.Cells.Item(55, 11).Value = 12000
$ws.Cells.Item(55, 12).Value = 12500
$ws.Cells.Item(55, 13).Value = 12250
$ws.Cells.Item(55, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(55, 16).Value = 204
$ws.Cells.Item(55, 17).Value = 60

$ws.Cells.Item(56, 4).Value = 44557
$ws.Cells.Item(56, 10).Value = 500
$ws.Cells.Item(56, 11).Value = 9500
$ws.Cells.Item(56, 12).Value = 10000
$ws.Cells.Item(56, 13).Value = 9750
$ws.Cells.Item(56, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(56, 16).Value = 162
$ws.Cells.Item(56, 17).Value = 60

$ws.Cells.Item(57, 4).Value = 44316
$ws.Cells.Item(57, 10).Value = 520
$ws.Cells.Item(57, 11).Value = 8000
$ws.Cells.Item(57, 12).Value = 9000
$ws.Cells.Item(57, 13).Value = 8500
$ws.Cells.Item(57, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(57, 16).Value = 142
$ws.Cells.Item(57, 17).Value = 60

$ws.Cells.Item(58, 4).Value = 44397
$ws.Cells.Item(58, 10).Value = 560
$ws.Cells.Item(58, 11).Value = 12000
$ws.Cells.Item(58, 12).Value = 12500
$ws.Cells.Item(58, 13).Value = 12250
$ws.Cells.Item(58, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(58, 16).Value = 204
$ws.Cells.Item(58, 17).Value = 60

$ws.Cells.Item(59, 4).Value = 44587
$ws.Cells.Item(59, 10).Value = 520
$ws.Cells.Item(59, 11).Value = 8000
$ws.Cells.Item(59, 12).Value = 9000
$ws.Cells.Item(59, 13).Value = 8500
$ws.Cells.Item(59, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(59, 16).Value = 170
$ws.Cells.Item(59, 17).Value = 50

$ws.Cells.Item(60, 4).Value = 44484
$ws.Cells.Item(60, 10).Value = 600
$ws.Cells.Item(60, 11).Value = 9000
$ws.Cells.Item(60, 12).Value = 10000
$ws.Cells.Item(60, 13).Value = 9500
$ws.Cells.Item(60, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(60, 16).Value = 158
$ws.Cells.Item(60, 17).Value = 60

$ws.Cells.Item(61, 4).Value = 44382
$ws.Cells.Item(61, 10).Value = 560
$ws.Cells.Item(61, 11).Value = 12000
$ws.Cells.Item(61, 12).Value = 13000
$ws.Cells.Item(61, 13).Value = 12500
$ws.Cells.Item(61, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(61, 16).Value = 208
$ws.Cells.Item(61, 17).Value = 60

$ws.Cells.Item(62, 4).Value = 44407
$ws.Cells.Item(62, 10).Value = 600
$ws.Cells.Item(62, 11).Value = 12500
$ws.Cells.Item(62, 12).Value = 13000
$ws.Cells.Item(62, 13).Value = 12750
$ws.Cells.Item(62, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(62, 16).Value = 212
$ws.Cells.Item(62, 17).Value = 60

$ws.Cells.Item(63, 4).Value = 44664
$ws.Cells.Item(63, 10).Value = 520
$ws.Cells.Item(63, 11).Value = 8000
$ws.Cells.Item(63, 12).Value = 9000
$ws.Cells.Item(63, 13).Value = 8500
$ws.Cells.Item(63, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(63, 16).Value = 170
$ws.Cells.Item(63, 17).Value = 50

$ws.Cells.Item(64, 4).Value = 44643
$ws.Cells.Item(64, 10).Value = 560
$ws.Cells.Item(64, 11).Value = 8000
$ws.Cells.Item(64, 12).Value = 9000
$ws.Cells.Item(64, 13).Value = 8500
$ws.Cells.Item(64, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(64, 16).Value = 170
$ws.Cells.Item(64, 17).Value = 50

$ws.Cells.Item(65, 4).Value = 44656
$ws.Cells.Item(65, 10).Value = 400
$ws.Cells.Item(65, 11).Value = 8000
$ws.Cells.Item(65, 12).Value = 9000
$ws.Cells.Item(65, 13).Value = 8500
$ws.Cells.Item(65, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(65, 16).Value = 170
$ws.Cells.Item(65, 17).Value = 50

$ws.Cells.Item(66, 4).Value = 44300
$ws.Cells.Item(66, 10).Value = 600
$ws.Cells.Item(66, 11).Value = 8000
$ws.Cells.Item(66, 12).Value = 8500
$ws.Cells.Item(66, 13).Value = 8250
$ws.Cells.Item(66, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(66, 16).Value = 138
$ws.Cells.Item(66, 17).Value = 60

$ws.Cells.Item(67, 4).Value = 44445
$ws.Cells.Item(67, 10).Value = 560
$ws.Cells.Item(67, 11).Value = 10000
$ws.Cells.Item(67, 12).Value = 11000
$ws.Cells.Item(67, 13).Value = 10500
$ws.Cells.Item(67, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(67, 16).Value = 175
$ws.Cells.Item(67, 17).Value = 60

$ws.Cells.Item(68, 4).Value = 44516
$ws.Cells.Item(68, 10).Value = 400
$ws.Cells.Item(68, 11).Value = 8000
$ws.Cells.Item(68, 12).Value = 9000
$ws.Cells.Item(68, 13).Value = 8500
$ws.Cells.Item(68, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(68, 16).Value = 142
$ws.Cells.Item(68, 17).Value = 60

$ws.Cells.Item(69, 4).Value = 44239
$ws.Cells.Item(69, 10).Value = 800
$ws.Cells.Item(69, 11).Value = 8000
$ws.Cells.Item(69, 12).Value = 9000
$ws.Cells.Item(69, 13).Value = 8500
$ws.Cells.Item(69, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(69, 16).Value = 142
$ws.Cells.Item(69, 17).Value = 60

$ws.Cells.Item(70, 4).Value = 44559
$ws.Cells.Item(70, 10).Value = 540
$ws.Cells.Item(70, 11).Value = 9500
$ws.Cells.Item(70, 12).Value = 10000
$ws.Cells.Item(70, 13).Value = 9750
$ws.Cells.Item(70, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(70, 16).Value = 162
$ws.Cells.Item(70, 17).Value = 60

$ws.Cells.Item(71, 4).Value = 44329
$ws.Cells.Item(71, 10).Value = 460
$ws.Cells.Item(71, 11).Value = 12000
$ws.Cells.Item(71, 12).Value = 13000
$ws.Cells.Item(71, 13).Value = 12500
$ws.Cells.Item(71, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(71, 16).Value = 208
$ws.Cells.Item(71, 17).Value = 60

$ws.Cells.Item(72, 4).Value = 44323
$ws.Cells.Item(72, 10).Value = 500
$ws.Cells.Item(72, 11).Value = 8000
$ws.Cells.Item(72, 12).Value = 9000
$ws.Cells.Item(72, 13).Value = 8500
$ws.Cells.Item(72, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(72, 16).Value = 142
$ws.Cells.Item(72, 17).Value = 60

$ws.Cells.Item(73, 4).Value = 44421
$ws.Cells.Item(73, 10).Value = 600
$ws.Cells.Item(73, 11).Value = 12000
$ws.Cells.Item(73, 12).Value = 12500
$ws.Cells.Item(73, 13).Value = 12250
$ws.Cells.Item(73, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(73, 16).Value = 204
$ws.Cells.Item(73, 17).Value = 60

$ws.Cells.Item(74, 4).Value = 44267
$ws.Cells.Item(74, 10).Value = 600
$ws.Cells.Item(74, 11).Value = 8000
$ws.Cells.Item(74, 12).Value = 8500
$ws.Cells.Item(74, 13).Value = 8250
$ws.Cells.Item(74, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(74, 16).Value = 138
$ws.Cells.Item(74, 17).Value = 60

$ws.Cells.Item(75, 4).Value = 44334
$ws.Cells.Item(75, 10).Value = 540
$ws.Cells.Item(75, 11).Value = 12500
$ws.Cells.Item(75, 12).Value = 13000
$ws.Cells.Item(75, 13).Value = 12750
$ws.Cells.Item(75, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(75, 16).Value = 212
$ws.Cells.Item(75, 17).Value = 60

$ws.Cells.Item(76, 4).Value = 44475
$ws.Cells.Item(76, 10).Value = 600
$ws.Cells.Item(76, 11).Value = 6000
$ws.Cells.Item(76, 12).Value = 7000
$ws.Cells.Item(76, 13).Value = 6500
$ws.Cells.Item(76, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(76, 16).Value = 108
$ws.Cells.Item(76, 17).Value = 60

$ws.Cells.Item(77, 4).Value = 44554
$ws.Cells.Item(77, 10).Value = 500
$ws.Cells.Item(77, 11).Value = 10000
$ws.Cells.Item(77, 12).Value = 11000
$ws.Cells.Item(77, 13).Value = 10500
$ws.Cells.Item(77, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(77, 16).Value = 175
$ws.Cells.Item(77, 17).Value = 60

$ws.Cells.Item(78, 4).Value = 44333
$ws.Cells.Item(78, 10).Value = 500
$ws.Cells.Item(78, 11).Value = 12000
$ws.Cells.Item(78, 12).Value = 13000
$ws.Cells.Item(78, 13).Value = 12500
$ws.Cells.Item(78, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(78, 16).Value = 208
$ws.Cells.Item(78, 17).Value = 60

$ws.Cells.Item(79, 4).Value = 44649
$ws.Cells.Item(79, 10).Value = 400
$ws.Cells.Item(79, 11).Value = 8000
$ws.Cells.Item(79, 12).Value = 9000
$ws.Cells.Item(79, 13).Value = 8500
$ws.Cells.Item(79, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(79, 16).Value = 170
$ws.Cells.Item(79, 17).Value = 50

$ws.Cells.Item(80, 4).Value = 44491
$ws.Cells.Item(80, 10).Value = 600
$ws.Cells.Item(80, 11).Value = 8500
$ws.Cells.Item(80, 12).Value = 9000
$ws.Cells.Item(80, 13).Value = 8750
$ws.Cells.Item(80, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(80, 16).Value = 146
$ws.Cells.Item(80, 17).Value = 60

$ws.Cells.Item(81, 4).Value = 44356
$ws.Cells.Item(81, 10).Value = 600
$ws.Cells.Item(81, 11).Value = 12000
$ws.Cells.Item(81, 12).Value = 13000
$ws.Cells.Item(81, 13).Value = 12500
$ws.Cells.Item(81, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(81, 16).Value = 208
$ws.Cells.Item(81, 17).Value = 60

$ws.Cells.Item(82, 4).Value = 44508
$ws.Cells.Item(82, 10).Value = 520
$ws.Cells.Item(82, 11).Value = 8000
$ws.Cells.Item(82, 12).Value = 9000
$ws.Cells.Item(82, 13).Value = 8500
$ws.Cells.Item(82, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(82, 16).Value = 142
$ws.Cells.Item(82, 17).Value = 60

$ws.Cells.Item(83, 4).Value = 44403
$ws.Cells.Item(83, 10).Value = 760
$ws.Cells.Item(83, 11).Value = 12000
$ws.Cells.Item(83, 12).Value = 12500
$ws.Cells.Item(83, 13).Value = 12250
$ws.Cells.Item(83, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(83, 16).Value = 204
$ws.Cells.Item(83, 17).Value = 60

$ws.Cells.Item(84, 4).Value = 44607
$ws.Cells.Item(84, 10).Value = 480
$ws.Cells.Item(84, 11).Value = 8500
$ws.Cells.Item(84, 12).Value = 9000
$ws.Cells.Item(84, 13).Value = 8750
$ws.Cells.Item(84, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(84, 16).Value = 175
$ws.Cells.Item(84, 17).Value = 50

$ws.Cells.Item(85, 4).Value = 44293
$ws.Cells.Item(85, 10).Value = 600
$ws.Cells.Item(85, 11).Value = 8000
$ws.Cells.Item(85, 12).Value = 9000
$ws.Cells.Item(85, 13).Value = 8500
$ws.Cells.Item(85, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(85, 16).Value = 142
$ws.Cells.Item(85, 17).Value = 60

$ws.Cells.Item(86, 4).Value = 44160
$ws.Cells.Item(86, 10).Value = 700
$ws.Cells.Item(86, 11).Value = 9000
$ws.Cells.Item(86, 12).Value = 10000
$ws.Cells.Item(86, 13).Value = 9500
$ws.Cells.Item(86, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(86, 16).Value = 158
$ws.Cells.Item(86, 17).Value = 60

$ws.Cells.Item(87, 4).Value = 44246
$ws.Cells.Item(87, 10).Value = 800
$ws.Cells.Item(87, 11).Value = 8000
$ws.Cells.Item(87, 12).Value = 8500
$ws.Cells.Item(87, 13).Value = 8250
$ws.Cells.Item(87, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(87, 16).Value = 138
$ws.Cells.Item(87, 17).Value = 60

$ws.Cells.Item(88, 4).Value = 44628
$ws.Cells.Item(88, 10).Value = 520
$ws.Cells.Item(88, 11).Value = 8500
$ws.Cells.Item(88, 12).Value = 9000
$ws.Cells.Item(88, 13).Value = 8750
$ws.Cells.Item(88, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(88, 16).Value = 175
$ws.Cells.Item(88, 17).Value = 50

$ws.Cells.Item(89, 4).Value = 44405
$ws.Cells.Item(89, 10).Value = 600
$ws.Cells.Item(89, 11).Value = 12000
$ws.Cells.Item(89, 12).Value = 12500
$ws.Cells.Item(89, 13).Value = 12250
$ws.Cells.Item(89, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(89, 16).Value = 204
$ws.Cells.Item(89, 17).Value = 60

$ws.Cells.Item(90, 4).Value = 44312
$ws.Cells.Item(90, 10).Value = 600
$ws.Cells.Item(90, 11).Value = 8000
$ws.Cells.Item(90, 12).Value = 9000
$ws.Cells.Item(90, 13).Value = 8500
$ws.Cells.Item(90, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(90, 16).Value = 142
$ws.Cells.Item(90, 17).Value = 60

$ws.Cells.Item(91, 4).Value = 44412
$ws.Cells.Item(91, 10).Value = 700
$ws.Cells.Item(91, 11).Value = 12500
$ws.Cells.Item(91, 12).Value = 13000
$ws.Cells.Item(91, 13).Value = 12750
$ws.Cells.Item(91, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(91, 16).Value = 212
$ws.Cells.Item(91, 17).Value = 60

$ws.Cells.Item(92, 4).Value = 44365
$ws.Cells.Item(92, 10).Value = 520
$ws.Cells.Item(92, 11).Value = 13000
$ws.Cells.Item(92, 12).Value = 14000
$ws.Cells.Item(92, 13).Value = 13500
$ws.Cells.Item(92, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(92, 16).Value = 225
$ws.Cells.Item(92, 17).Value = 60

$ws.Cells.Item(93, 4).Value = 44586
$ws.Cells.Item(93, 10).Value = 600
$ws.Cells.Item(93, 11).Value = 8000
$ws.Cells.Item(93, 12).Value = 9000
$ws.Cells.Item(93, 13).Value = 8500
$ws.Cells.Item(93, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(93, 16).Value = 170
$ws.Cells.Item(93, 17).Value = 50

$ws.Cells.Item(94, 4).Value = 44468
$ws.Cells.Item(94, 10).Value = 600
$ws.Cells.Item(94, 11).Value = 8000
$ws.Cells.Item(94, 12).Value = 9000
$ws.Cells.Item(94, 13).Value = 8500
$ws.Cells.Item(94, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(94, 16).Value = 142
$ws.Cells.Item(94, 17).Value = 60

$ws.Cells.Item(95, 4).Value = 44370
$ws.Cells.Item(95, 10).Value = 600
$ws.Cells.Item(95, 11).Value = 13000
$ws.Cells.Item(95, 12).Value = 14000
$ws.Cells.Item(95, 13).Value = 13500
$ws.Cells.Item(95, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(95, 16).Value = 225
$ws.Cells.Item(95, 17).Value = 60

$ws.Cells.Item(96, 4).Value = 44265
$ws.Cells.Item(96, 10).Value = 720
$ws.Cells.Item(96, 11).Value = 8000
$ws.Cells.Item(96, 12).Value = 9000
$ws.Cells.Item(96, 13).Value = 8500
$ws.Cells.Item(96, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(96, 16).Value = 142
$ws.Cells.Item(96, 17).Value = 60

$ws.Cells.Item(97, 4).Value = 44463
$ws.Cells.Item(97, 10).Value = 600
$ws.Cells.Item(97, 11).Value = 9000
$ws.Cells.Item(97, 12).Value = 10000
$ws.Cells.Item(97, 13).Value = 9500
$ws.Cells.Item(97, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(97, 16).Value = 158
$ws.Cells.Item(97, 17).Value = 60

$ws.Cells.Item(98, 4).Value = 44544
$ws.Cells.Item(98, 10).Value = 600
$ws.Cells.Item(98, 11).Value = 10000
$ws.Cells.Item(98, 12).Value = 11000
$ws.Cells.Item(98, 13).Value = 10500
$ws.Cells.Item(98, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(98, 16).Value = 175
$ws.Cells.Item(98, 17).Value = 60

$ws.Cells.Item(99, 4).Value = 44342
$ws.Cells.Item(99, 10).Value = 600
$ws.Cells.Item(99, 11).Value = 12000
$ws.Cells.Item(99, 12).Value = 13000
$ws.Cells.Item(99, 13).Value = 12500
$ws.Cells.Item(99, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(99, 16).Value = 208
$ws.Cells.Item(99, 17).Value = 60

$ws.Cells.Item(100, 4).Value = 44568
$ws.Cells.Item(100, 10).Value = 700
$ws.Cells.Item(100, 11).Value = 8000
$ws.Cells.Item(100, 12).Value = 9000
$ws.Cells.Item(100, 13).Value = 8500
$ws.Cells.Item(100, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(100, 16).Value = 142
$ws.Cells.Item(100, 17).Value = 60

$ws.Cells.Item(101, 4).Value = 44162
$ws.Cells.Item(101, 10).Value = 600
$ws.Cells.Item(101, 11).Value = 9000
$ws.Cells.Item(101, 12).Value = 10000
$ws.Cells.Item(101, 13).Value = 9500
$ws.Cells.Item(101, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(101, 16).Value = 158
$ws.Cells.Item(101, 17).Value = 60

$ws.Cells.Item(102, 4).Value = 44512
$ws.Cells.Item(102, 10).Value = 600
$ws.Cells.Item(102, 11).Value = 8000
$ws.Cells.Item(102, 12).Value = 9000
$ws.Cells.Item(102, 13).Value = 8500
$ws.Cells.Item(102, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(102, 16).Value = 142
$ws.Cells.Item(102, 17).Value = 60

$ws.Cells.Item(103, 4).Value = 44582
$ws.Cells.Item(103, 10).Value = 600
$ws.Cells.Item(103, 11).Value = 8500
$ws.Cells.Item(103, 12).Value = 9000
$ws.Cells.Item(103, 13).Value = 8750
$ws.Cells.Item(103, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(103, 16).Value = 175
$ws.Cells.Item(103, 17).Value = 50

$ws.Cells.Item(104, 4).Value = 44251
$ws.Cells.Item(104, 10).Value = 600
$ws.Cells.Item(104, 11).Value = 8000
$ws.Cells.Item(104, 12).Value = 9000
$ws.Cells.Item(104, 13).Value = 8500
$ws.Cells.Item(104, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(104, 16).Value = 142
$ws.Cells.Item(104, 17).Value = 60

$ws.Cells.Item(105, 4).Value = 44258
$ws.Cells.Item(105, 10).Value = 700
$ws.Cells.Item(105, 11).Value = 8000
$ws.Cells.Item(105, 12).Value = 8500
$ws.Cells.Item(105, 13).Value = 8250
$ws.Cells.Item(105, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(105, 16).Value = 138
$ws.Cells.Item(105, 17).Value = 60

$ws.Cells.Item(106, 4).Value = 44372
$ws.Cells.Item(106, 10).Value = 560
$ws.Cells.Item(106, 11).Value = 13000
$ws.Cells.Item(106, 12).Value = 14000
$ws.Cells.Item(106, 13).Value = 13500
$ws.Cells.Item(106, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(106, 16).Value = 225
$ws.Cells.Item(106, 17).Value = 60

$ws.Cells.Item(107, 4).Value = 44524
$ws.Cells.Item(107, 10).Value = 540
$ws.Cells.Item(107, 11).Value = 9000
$ws.Cells.Item(107, 12).Value = 10000
$ws.Cells.Item(107, 13).Value = 9500
$ws.Cells.Item(107, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(107, 16).Value = 158
$ws.Cells.Item(107, 17).Value = 60

$ws.Cells.Item(108, 4).Value = 44452
$ws.Cells.Item(108, 10).Value = 560
$ws.Cells.Item(108, 11).Value = 9000
$ws.Cells.Item(108, 12).Value = 10000
$ws.Cells.Item(108, 13).Value = 9500
$ws.Cells.Item(108, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(108, 16).Value = 158
$ws.Cells.Item(108, 17).Value = 60

$ws.Cells.Item(109, 4).Value = 44634
$ws.Cells.Item(109, 10).Value = 480
$ws.Cells.Item(109, 11).Value = 8500
$ws.Cells.Item(109, 12).Value = 9000
$ws.Cells.Item(109, 13).Value = 8750
$ws.Cells.Item(109, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(109, 16).Value = 175
$ws.Cells.Item(109, 17).Value = 50

$ws.Cells.Item(110, 4).Value = 44385
$ws.Cells.Item(110, 10).Value = 560
$ws.Cells.Item(110, 11).Value = 11000
$ws.Cells.Item(110, 12).Value = 12000
$ws.Cells.Item(110, 13).Value = 11500
$ws.Cells.Item(110, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(110, 16).Value = 192
$ws.Cells.Item(110, 17).Value = 60

$ws.Cells.Item(111, 4).Value = 44552
$ws.Cells.Item(111, 10).Value = 540
$ws.Cells.Item(111, 11).Value = 10500
$ws.Cells.Item(111, 12).Value = 11000
$ws.Cells.Item(111, 13).Value = 10750
$ws.Cells.Item(111, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(111, 16).Value = 179
$ws.Cells.Item(111, 17).Value = 60

$ws.Cells.Item(112, 4).Value = 44498
$ws.Cells.Item(112, 10).Value = 560
$ws.Cells.Item(112, 11).Value = 8000
$ws.Cells.Item(112, 12).Value = 9000
$ws.Cells.Item(112, 13).Value = 8500
$ws.Cells.Item(112, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(112, 16).Value = 142
$ws.Cells.Item(112, 17).Value = 60

$ws.Cells.Item(113, 4).Value = 44620
$ws.Cells.Item(113, 10).Value = 520
$ws.Cells.Item(113, 11).Value = 8000
$ws.Cells.Item(113, 12).Value = 9000
$ws.Cells.Item(113, 13).Value = 8500
$ws.Cells.Item(113, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(113, 16).Value = 170
$ws.Cells.Item(113, 17).Value = 50

$ws.Cells.Item(114, 4).Value = 44371
$ws.Cells.Item(114, 10).Value = 560
$ws.Cells.Item(114, 11).Value = 13000
$ws.Cells.Item(114, 12).Value = 14000
$ws.Cells.Item(114, 13).Value = 13500
$ws.Cells.Item(114, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(114, 16).Value = 225
$ws.Cells.Item(114, 17).Value = 60

$ws.Cells.Item(115, 4).Value = 44600
$ws.Cells.Item(115, 10).Value = 520
$ws.Cells.Item(115, 11).Value = 8500
$ws.Cells.Item(115, 12).Value = 9000
$ws.Cells.Item(115, 13).Value = 8750
$ws.Cells.Item(115, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(115, 16).Value = 175
$ws.Cells.Item(115, 17).Value = 50

$ws.Cells.Item(116, 4).Value = 44237
$ws.Cells.Item(116, 10).Value = 600
$ws.Cells.Item(116, 11).Value = 8000
$ws.Cells.Item(116, 12).Value = 9000
$ws.Cells.Item(116, 13).Value = 8500
$ws.Cells.Item(116, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(116, 16).Value = 142
$ws.Cells.Item(116, 17).Value = 60

$ws.Cells.Item(117, 4).Value = 44326
$ws.Cells.Item(117, 10).Value = 500
$ws.Cells.Item(117, 11).Value = 9000
$ws.Cells.Item(117, 12).Value = 10000
$ws.Cells.Item(117, 13).Value = 9500
$ws.Cells.Item(117, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(117, 16).Value = 158
$ws.Cells.Item(117, 17).Value = 60

$ws.Cells.Item(118, 4).Value = 44473
$ws.Cells.Item(118, 10).Value = 700
$ws.Cells.Item(118, 11).Value = 7000
$ws.Cells.Item(118, 12).Value = 8000
$ws.Cells.Item(118, 13).Value = 7500
$ws.Cells.Item(118, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(118, 16).Value = 125
$ws.Cells.Item(118, 17).Value = 60

$ws.Cells.Item(119, 4).Value = 44428
$ws.Cells.Item(119, 10).Value = 600
$ws.Cells.Item(119, 11).Value = 12000
$ws.Cells.Item(119, 12).Value = 13000
$ws.Cells.Item(119, 13).Value = 12500
$ws.Cells.Item(119, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(119, 16).Value = 208
$ws.Cells.Item(119, 17).Value = 60

$ws.Cells.Item(120, 4).Value = 44406
$ws.Cells.Item(120, 10).Value = 600
$ws.Cells.Item(120, 11).Value = 12000
$ws.Cells.Item(120, 12).Value = 12500
$ws.Cells.Item(120, 13).Value = 12250
$ws.Cells.Item(120, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(120, 16).Value = 204
$ws.Cells.Item(120, 17).Value = 60

$ws.Cells.Item(121, 4).Value = 44386
$ws.Cells.Item(121, 10).Value = 560
$ws.Cells.Item(121, 11).Value = 11000
$ws.Cells.Item(121, 12).Value = 12000
$ws.Cells.Item(121, 13).Value = 11500
$ws.Cells.Item(121, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(121, 16).Value = 192
$ws.Cells.Item(121, 17).Value = 60

$ws.Cells.Item(122, 4).Value = 44427
$ws.Cells.Item(122, 10).Value = 560
$ws.Cells.Item(122, 11).Value = 12000
$ws.Cells.Item(122, 12).Value = 13000
$ws.Cells.Item(122, 13).Value = 12500
$ws.Cells.Item(122, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(122, 16).Value = 208
$ws.Cells.Item(122, 17).Value = 60

$ws.Cells.Item(123, 4).Value = 44389
$ws.Cells.Item(123, 10).Value = 660
$ws.Cells.Item(123, 11).Value = 11500
$ws.Cells.Item(123, 12).Value = 12000
$ws.Cells.Item(123, 13).Value = 11750
$ws.Cells.Item(123, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(123, 16).Value = 196
$ws.Cells.Item(123, 17).Value = 60

$ws.Cells.Item(124, 4).Value = 44641
$ws.Cells.Item(124, 10).Value = 500
$ws.Cells.Item(124, 11).Value = 8500
$ws.Cells.Item(124, 12).Value = 9000
$ws.Cells.Item(124, 13).Value = 8750
$ws.Cells.Item(124, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(124, 16).Value = 175
$ws.Cells.Item(124, 17).Value = 50

$ws.Cells.Item(125, 4).Value = 44426
$ws.Cells.Item(125, 10).Value = 600
$ws.Cells.Item(125, 11).Value = 12500
$ws.Cells.Item(125, 12).Value = 13000
$ws.Cells.Item(125, 13).Value = 12750
$ws.Cells.Item(125, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(125, 16).Value = 212
$ws.Cells.Item(125, 17).Value = 60

$ws.Cells.Item(126, 4).Value = 44335
$ws.Cells.Item(126, 10).Value = 600
$ws.Cells.Item(126, 11).Value = 12000
$ws.Cells.Item(126, 12).Value = 13000
$ws.Cells.Item(126, 13).Value = 12500
$ws.Cells.Item(126, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(126, 16).Value = 208
$ws.Cells.Item(126, 17).Value = 60

$ws.Cells.Item(127, 4).Value = 44330
$ws.Cells.Item(127, 10).Value = 520
$ws.Cells.Item(127, 11).Value = 12000
$ws.Cells.Item(127, 12).Value = 13000
$ws.Cells.Item(127, 13).Value = 12500
$ws.Cells.Item(127, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(127, 16).Value = 208
$ws.Cells.Item(127, 17).Value = 60

$ws.Cells.Item(128, 4).Value = 44343
$ws.Cells.Item(128, 10).Value = 500
$ws.Cells.Item(128, 11).Value = 12000
$ws.Cells.Item(128, 12).Value = 13000
$ws.Cells.Item(128, 13).Value = 12500
$ws.Cells.Item(128, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(128, 16).Value = 208
$ws.Cells.Item(128, 17).Value = 60

$ws.Cells.Item(129, 4).Value = 44358
$ws.Cells.Item(129, 10).Value = 540
$ws.Cells.Item(129, 11).Value = 11500
$ws.Cells.Item(129, 12).Value = 12000
$ws.Cells.Item(129, 13).Value = 11750
$ws.Cells.Item(129, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(129, 16).Value = 196
$ws.Cells.Item(129, 17).Value = 60

$ws.Cells.Item(130, 4).Value = 44349
$ws.Cells.Item(130, 10).Value = 600
$ws.Cells.Item(130, 11).Value = 12000
$ws.Cells.Item(130, 12).Value = 12500
$ws.Cells.Item(130, 13).Value = 12250
$ws.Cells.Item(130, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(130, 16).Value = 204
$ws.Cells.Item(130, 17).Value = 60

$ws.Cells.Item(131, 4).Value = 44466
$ws.Cells.Item(131, 10).Value = 600
$ws.Cells.Item(131, 11).Value = 8500
$ws.Cells.Item(131, 12).Value = 9000
$ws.Cells.Item(131, 13).Value = 8750
$ws.Cells.Item(131, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(131, 16).Value = 146
$ws.Cells.Item(131, 17).Value = 60

$ws.Cells.Item(132, 4).Value = 44307
$ws.Cells.Item(132, 10).Value = 600
$ws.Cells.Item(132, 11).Value = 8000
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = 8500
$ws.Cells.Item(132, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(132, 16).Value = 142
$ws.Cells.Item(132, 17).Value = 60

